$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (rich-text shared strings) -- "Volume 29 Number 44"
#    becomes "...Number 45", and the report-week dates move forward one week.
# ---------------------------------------------------------------------------

# A8 = "Volume 29   Number  44"  ->  "...  45"
$a8 = $ws.Range("A8")
$a8Text = $a8.Text
$issueStart = $a8Text.Length - 1
$issueLen = 2
$a8.Characters($issueStart, $issueLen).Text = "45"

# C9 = "Report Covering the Week  10/31/2022  Through  11/6/2022"
#   -> "Report Covering the Week  11/7/2022  Through  11/13/2022"
# Replace the later run first so the earlier run's character offset stays valid.
$c9 = $ws.Range("C9")
$c9Text = $c9.Text
$throughStart = $c9Text.Length - "11/6/2022".Length + 1
$c9.Characters($throughStart, "11/6/2022".Length).Text = "11/13/2022"

$c9Text2 = $c9.Text
$weekStart = "Report Covering the Week  ".Length + 1
$c9.Characters($weekStart, "10/31/2022".Length).Text = "11/7/2022"

# ---------------------------------------------------------------------------
# 2) Weekly crime statistics grid (rows 14-30) -- refreshed counts/percentages
# ---------------------------------------------------------------------------

$ws.Range("L14").Value = -52.941176470588
$ws.Range("N14").Value = -61.904761904761

$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 32
$ws.Range("J15").Value = 41
$ws.Range("K15").Value = -21.951219512195
$ws.Range("L15").Value = -31.914893617021
$ws.Range("M15").Value = -37.254901960784
$ws.Range("N15").Value = -54.929577464788

$ws.Range("C16").Value = 5
$ws.Range("E16").Value = -16.666666666666
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = 26.315789473684
$ws.Range("I16").Value = 235
$ws.Range("J16").Value = 178
$ws.Range("K16").Value = 32.022471910112
$ws.Range("L16").Value = 14.634146341463
$ws.Range("M16").Value = -36.486486486486
$ws.Range("N16").Value = -80.562448304383

$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -31.25
$ws.Range("F17").Value = 57
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 23.913043478260
$ws.Range("I17").Value = 684
$ws.Range("J17").Value = 526
$ws.Range("K17").Value = 30.038022813688
$ws.Range("L17").Value = 28.813559322033
$ws.Range("M17").Value = 58.333333333333
$ws.Range("N17").Value = -35.955056179775

$ws.Range("C18").Value = 9
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 275
$ws.Range("J18").Value = 213
$ws.Range("K18").Value = 29.107981220657
$ws.Range("L18").Value = -3.508771929824
$ws.Range("M18").Value = -49.908925318761
$ws.Range("N18").Value = -90.833333333333

$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 31
$ws.Range("E19").Value = -29.032258064516
$ws.Range("F19").Value = 131
$ws.Range("G19").Value = 94
$ws.Range("H19").Value = 39.361702127659
$ws.Range("I19").Value = 1192
$ws.Range("J19").Value = 943
$ws.Range("K19").Value = 26.405090137857
$ws.Range("L19").Value = 42.413381123058
$ws.Range("M19").Value = 40.898345153664
$ws.Range("N19").Value = -17.049408489909

$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 68
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = 119.354838709677
$ws.Range("I20").Value = 432
$ws.Range("J20").Value = 228
$ws.Range("K20").Value = 89.473684210526
$ws.Range("L20").Value = 105.714285714286
$ws.Range("M20").Value = 39.805825242718
$ws.Range("N20").Value = -89.828113962797

$ws.Range("C21").Value = 68
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = 4.615384615384
$ws.Range("F21").Value = 305
$ws.Range("G21").Value = 207
$ws.Range("H21").Value = 47.342995169082
$ws.Range("I21").Value = 2858
$ws.Range("J21").Value = 2141
$ws.Range("K21").Value = 33.489023820644
$ws.Range("L21").Value = 34.052532833020
$ws.Range("M21").Value = 11.162971606378
$ws.Range("N21").Value = -74.142766669682

# Row 23: D23/E23 flip from plain numbers to the "N/A" placeholder text used
# elsewhere in the sheet (shared strings "0" and "***.*"), re-using the
# formatting of existing placeholder cells so the style index matches.
$ws.Range("C23").Value = 1

$d23 = $ws.Range("D23")
$d23.Value = "'0"
$ws.Range("D14").Copy()
$d23.PasteSpecial(-4122)

$e23 = $ws.Range("E23")
$e23.Value = "***.*"
$ws.Range("E22").Copy()
$e23.PasteSpecial(-4122)

$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -58.333333333333
$ws.Range("I23").Value = 84
$ws.Range("K23").Value = -10.638297872340
$ws.Range("L23").Value = 27.272727272727
$ws.Range("M23").Value = 33.333333333333

$ws.Range("C24").Value = 94
$ws.Range("D24").Value = 66
$ws.Range("E24").Value = 42.424242424242
$ws.Range("F24").Value = 344
$ws.Range("G24").Value = 226
$ws.Range("H24").Value = 52.212389380531
$ws.Range("I24").Value = 3438
$ws.Range("J24").Value = 2173
$ws.Range("K24").Value = 58.214450069029
$ws.Range("L24").Value = 54.032258064516
$ws.Range("M24").Value = 4.213397999393

$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 24
$ws.Range("G25").Value = 116
$ws.Range("H25").Value = 2.586206896551
$ws.Range("I25").Value = 1440
$ws.Range("J25").Value = 1177
$ws.Range("K25").Value = 22.344944774851
$ws.Range("L25").Value = 32.474701011959
$ws.Range("M25").Value = -19.732441471571

$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 100
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 66
$ws.Range("J26").Value = 79
$ws.Range("K26").Value = -16.455696202531
$ws.Range("L26").Value = -13.157894736842

$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 66.666666666666
$ws.Range("G27").Value = 21
$ws.Range("H27").Value = -38.095238095238
$ws.Range("I27").Value = 153
$ws.Range("J27").Value = 131
$ws.Range("K27").Value = 16.793893129771
$ws.Range("L27").Value = 61.052631578947

# Row 28: D28/E28 flip from plain numbers to the "N/A" placeholder text.
$d28 = $ws.Range("D28")
$d28.Value = "'0"
$ws.Range("D14").Copy()
$d28.PasteSpecial(-4122)

$e28 = $ws.Range("E28")
$e28.Value = "***.*"
$ws.Range("E22").Copy()
$e28.PasteSpecial(-4122)

$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = -15.789473684210
$ws.Range("N28").Value = -68.627450980392

# Row 29: D29/E29 flip from plain numbers to the "N/A" placeholder text.
$d29 = $ws.Range("D29")
$d29.Value = "'0"
$ws.Range("D14").Copy()
$d29.PasteSpecial(-4122)

$e29 = $ws.Range("E29")
$e29.Value = "***.*"
$ws.Range("E22").Copy()
$e29.PasteSpecial(-4122)

$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("L29").Value = -23.529411764705
$ws.Range("N29").Value = -71.111111111111

# Row 30: D30/E30/F30 flip from the "N/A" placeholder text back to plain
# numbers, re-using the formatting of existing numeric cells so the style
# index matches (General "#,##0" for counts, percent format for the change).
$d30 = $ws.Range("D30")
$ws.Range("C23").Copy()
$d30.PasteSpecial(-4122)
$d30.Value = 1

$e30 = $ws.Range("E30")
$ws.Range("H23").Copy()
$e30.PasteSpecial(-4122)
$e30.Value = -100

$f30 = $ws.Range("F30")
$ws.Range("C23").Copy()
$f30.PasteSpecial(-4122)
$f30.Value = 1

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 22
$ws.Range("J30").Value = 17
$ws.Range("K30").Value = 29.411764705882
$ws.Range("L30").Value = 22.222222222222
